# Natmi following Dr Hou advice
# Rewrites the LR-pair data rows (2-13) on Sheet1 with updated values,
# expanding the table from 8 data rows (Sending cluster x Target cluster
# combos for ECs/FAPs) to 12 rows (adding the sCs sending-cluster group).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col5a3"
$ws.Range("C2").Value = "Sdc3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.4980406666666666
$ws.Range("H2").Value = 1.494122
$ws.Range("I2").Value = 0.004393808999309369
$ws.Range("J2").Value = 0.00439380899930937
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 24.44779533333333
$ws.Range("N2").Value = 73.343386
$ws.Range("O2").Value = 0.1890645123346783
$ws.Range("P2").Value = 0.1890645123346783
$ws.Range("Q2").Value = 12.17599628634355
$ws.Range("R2").Value = 109.583966577092
$ws.Range("S2").Value = 0.0008307133557461467
$ws.Range("T2").Value = 0.0008307133557461469

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col5a3"
$ws.Range("C3").Value = "Sdc3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.4980406666666666
$ws.Range("H3").Value = 1.494122
$ws.Range("I3").Value = 0.004393808999309369
$ws.Range("J3").Value = 0.00439380899930937
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 18.11074233333333
$ws.Range("N3").Value = 54.332227
$ws.Range("O3").Value = 0.1400575643155068
$ws.Range("P3").Value = 0.1400575643155068
$ws.Range("Q3").Value = 9.019886185521555
$ws.Range("R3").Value = 81.178975669694
$ws.Range("S3").Value = 0.0006153861865108247
$ws.Range("T3").Value = 0.0006153861865108248

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col5a3"
$ws.Range("C4").Value = "Sdc3"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4980406666666666
$ws.Range("H4").Value = 1.494122
$ws.Range("I4").Value = 0.004393808999309369
$ws.Range("J4").Value = 0.00439380899930937
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 80.684877
$ws.Range("N4").Value = 242.054631
$ws.Range("O4").Value = 0.6239682030546764
$ws.Range("P4").Value = 0.6239682030546765
$ws.Range("Q4").Value = 40.184349930998
$ws.Range("R4").Value = 361.659149378982
$ws.Range("S4").Value = 0.002741597105864533
$ws.Range("T4").Value = 0.002741597105864533

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col5a3"
$ws.Range("C5").Value = "Sdc3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4980406666666666
$ws.Range("H5").Value = 1.494122
$ws.Range("I5").Value = 0.004393808999309369
$ws.Range("J5").Value = 0.00439380899930937
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 6.065862
$ws.Range("N5").Value = 18.197586
$ws.Range("O5").Value = 0.0469097202951384
$ws.Range("P5").Value = 0.04690972029513841
$ws.Range("Q5").Value = 3.021045954388
$ws.Range("R5").Value = 27.189413589492
$ws.Range("S5").Value = 0.0002061123511878645
$ws.Range("T5").Value = 0.0002061123511878645

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col5a3"
$ws.Range("C6").Value = "Sdc3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 92.11319466666667
$ws.Range("H6").Value = 276.339584
$ws.Range("I6").Value = 0.8126400327714922
$ws.Range("J6").Value = 0.8126400327714922
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 24.44779533333333
$ws.Range("N6").Value = 73.343386
$ws.Range("O6").Value = 0.1890645123346783
$ws.Range("P6").Value = 0.1890645123346783
$ws.Range("Q6").Value = 2251.964530710158
$ws.Range("R6").Value = 20267.68077639142
$ws.Range("S6").Value = 0.1536413914995792
$ws.Range("T6").Value = 0.1536413914995792

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col5a3"
$ws.Range("C7").Value = "Sdc3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 92.11319466666667
$ws.Range("H7").Value = 276.339584
$ws.Range("I7").Value = 0.8126400327714922
$ws.Range("J7").Value = 0.8126400327714922
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 18.11074233333333
$ws.Range("N7").Value = 54.332227
$ws.Range("O7").Value = 0.1400575643155068
$ws.Range("P7").Value = 0.1400575643155068
$ws.Range("Q7").Value = 1668.238334108175
$ws.Range("R7").Value = 15014.14500697357
$ws.Range("S7").Value = 0.1138163836552489
$ws.Range("T7").Value = 0.1138163836552489

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col5a3"
$ws.Range("C8").Value = "Sdc3"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 92.11319466666667
$ws.Range("H8").Value = 276.339584
$ws.Range("I8").Value = 0.8126400327714922
$ws.Range("J8").Value = 0.8126400327714922
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 80.684877
$ws.Range("N8").Value = 242.054631
$ws.Range("O8").Value = 0.6239682030546764
$ws.Range("P8").Value = 0.6239682030546765
$ws.Range("Q8").Value = 7432.141781757056
$ws.Range("R8").Value = 66889.2760358135
$ws.Range("S8").Value = 0.5070615409787214
$ws.Range("T8").Value = 0.5070615409787214

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col5a3"
$ws.Range("C9").Value = "Sdc3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 92.11319466666667
$ws.Range("H9").Value = 276.339584
$ws.Range("I9").Value = 0.8126400327714922
$ws.Range("J9").Value = 0.8126400327714922
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 6.065862
$ws.Range("N9").Value = 18.197586
$ws.Range("O9").Value = 0.0469097202951384
$ws.Range("P9").Value = 0.04690972029513841
$ws.Range("Q9").Value = 558.745927227136
$ws.Range("R9").Value = 5028.713345044224
$ws.Range("S9").Value = 0.0381207166379428
$ws.Range("T9").Value = 0.03812071663794281

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Col5a3"
$ws.Range("C10").Value = "Sdc3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 20.73931466666667
$ws.Range("H10").Value = 62.217944
$ws.Range("I10").Value = 0.1829661582291984
$ws.Range("J10").Value = 0.1829661582291984
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 24.44779533333333
$ws.Range("N10").Value = 73.343386
$ws.Range("O10").Value = 0.1890645123346783
$ws.Range("P10").Value = 0.1890645123346783
$ws.Range("Q10").Value = 507.0305203242649
$ws.Range("R10").Value = 4563.274682918384
$ws.Range("S10").Value = 0.03459240747935299
$ws.Range("T10").Value = 0.03459240747935299

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Col5a3"
$ws.Range("C11").Value = "Sdc3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 20.73931466666667
$ws.Range("H11").Value = 62.217944
$ws.Range("I11").Value = 0.1829661582291984
$ws.Range("J11").Value = 0.1829661582291984
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 18.11074233333333
$ws.Range("N11").Value = 54.332227
$ws.Range("O11").Value = 0.1400575643155068
$ws.Range("P11").Value = 0.1400575643155068
$ws.Range("Q11").Value = 375.604384097921
$ws.Range("R11").Value = 3380.439456881288
$ws.Range("S11").Value = 0.02562579447374716
$ws.Range("T11").Value = 0.02562579447374716

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Col5a3"
$ws.Range("C12").Value = "Sdc3"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 20.73931466666667
$ws.Range("H12").Value = 62.217944
$ws.Range("I12").Value = 0.1829661582291984
$ws.Range("J12").Value = 0.1829661582291984
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 80.684877
$ws.Range("N12").Value = 242.054631
$ws.Range("O12").Value = 0.6239682030546764
$ws.Range("P12").Value = 0.6239682030546765
$ws.Range("Q12").Value = 1673.349052944296
$ws.Range("R12").Value = 15060.14147649866
$ws.Range("S12").Value = 0.1141650649700905
$ws.Range("T12").Value = 0.1141650649700905

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Col5a3"
$ws.Range("C13").Value = "Sdc3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 20.73931466666667
$ws.Range("H13").Value = 62.217944
$ws.Range("I13").Value = 0.1829661582291984
$ws.Range("J13").Value = 0.1829661582291984
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 6.065862
$ws.Range("N13").Value = 18.197586
$ws.Range("O13").Value = 0.0469097202951384
$ws.Range("P13").Value = 0.04690972029513841
$ws.Range("Q13").Value = 125.801820742576
$ws.Range("R13").Value = 1132.216386683184
$ws.Range("S13").Value = 0.008582891306007734
$ws.Range("T13").Value = 0.008582891306007735

